# Weekly update: add a new "Cebollín" price entry for the latest week.
# This inserts a new row at row 203 (pushing all subsequent rows down by
# one, e.g. old row 237 becomes row 238) and fills it with the new
# observation dated 45275 (2023-12-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 203; existing rows 203:237
# shift down to 204:238, preserving all of their values/styles.
$ws.Rows("203:203").Insert()

# Populate the newly inserted row 203 with the new weekly observation.
$ws.Cells.Item(203, 1).Value = 7
$ws.Cells.Item(203, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(203, 3).Value = "Ñuble"
$ws.Cells.Item(203, 4).Value = 45275
$ws.Cells.Item(203, 5).Value = 16
$ws.Cells.Item(203, 6).Value = 100112037
$ws.Cells.Item(203, 7).Value = "Cebollín"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 150
$ws.Cells.Item(203, 11).Value = 6000
$ws.Cells.Item(203, 12).Value = 6000
$ws.Cells.Item(203, 13).Value = 6000
$ws.Cells.Item(203, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(203, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(203, 16).Value = 167
$ws.Cells.Item(203, 17).Value = 36
$ws.Cells.Item(203, 18).Value = "Hortaliza"
